$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of data was added for "Jengibre" - Terminal La Palmera de La Serena.
# It belongs chronologically right before the current row 155, so insert a new
# row there (shifting the existing rows 155:187 down to 156:188) and then fill
# it in with the new record's values. All other previously-shifted rows keep
# their original data automatically.
$ws.Rows("155:155").Insert()

$newRow = 155
$ws.Cells.Item($newRow, 1).Value = 8
$ws.Cells.Item($newRow, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 45209
$ws.Cells.Item($newRow, 5).Value = 4
$ws.Cells.Item($newRow, 6).Value = 100114007
$ws.Cells.Item($newRow, 7).Value = "Jengibre"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 500
$ws.Cells.Item($newRow, 11).Value = 21000
$ws.Cells.Item($newRow, 12).Value = 22000
$ws.Cells.Item($newRow, 13).Value = 21500
$ws.Cells.Item($newRow, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item($newRow, 15).Value = "Perú"
$ws.Cells.Item($newRow, 16).Value = 1654
$ws.Cells.Item($newRow, 17).Value = 13
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
